$wb = $excel.ActiveWorkbook

# --- Sheets (by tab order): 1 = salsa, 2 = son, 3 = rumba ---
$salsa = $wb.Worksheets.Item(1)
$son   = $wb.Worksheets.Item(2)
$rumba = $wb.Worksheets.Item(3)

# --- salsa (sheet1): add a new row 8 with a new URI string, which also
#     grows the sheet dimension from A1:D6 to A1:D8 ---
$salsa.Range("A8").Value = "3Pacy6CMa8HPNVfeA3wkPQ,"

# --- son (sheet2): selection moves from A8 to A7; it stops being the
#     active/selected tab (rumba takes over below) ---
[void]$son.Range("A7").Select()

# --- salsa selection ends on A6 ---
[void]$salsa.Range("A6").Select()

# --- rumba (sheet3) becomes the active / selected tab, selection stays A7 ---
[void]$rumba.Range("A7").Select()
